$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-23 Tuesday" "2024-01-24 Wednesday"

Replace-Text "816÷7=" "686÷2="
Replace-Text "141÷9=" "278÷7="
Replace-Text "840÷3=" "660÷4="
Replace-Text "599÷4=" "279÷4="
Replace-Text "421÷7=" "399÷9="

Replace-Text "396÷7=" "401÷9="
Replace-Text "224÷5=" "496÷7="
Replace-Text "364÷6=" "694÷7="
Replace-Text "530÷3=" "430÷7="
Replace-Text "128÷8=" "524÷8="

Replace-Text "182÷9=" "169÷2="
Replace-Text "307÷7=" "737÷9="
Replace-Text "353÷3=" "186÷4="
Replace-Text "316÷5=" "430÷2="
Replace-Text "478÷9=" "314÷9="

Replace-Text "542÷8=" "847÷7="
Replace-Text "253÷6=" "808÷4="
Replace-Text "498÷6=" "318÷6="
Replace-Text "647÷2=" "639÷2="
Replace-Text "634÷6=" "976÷7="

Replace-Text "489÷5=" "568÷5="
Replace-Text "950÷2=" "622÷2="
Replace-Text "302÷5=" "833÷8="
Replace-Text "351÷3=" "630÷4="
Replace-Text "901÷2=" "692÷5="
